$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'ORGANIZATION_TYPE'
$ws.Range("B2").Value = 1333
$ws.Range("A3").Value = 'EXT_SOURCE_1'
$ws.Range("B3").Value = 781
$ws.Range("A4").Value = 'EXT_SOURCE_3'
$ws.Range("B4").Value = 736
$ws.Range("A5").Value = 'EXT_SOURCE_2'
$ws.Range("B5").Value = 700
$ws.Range("A6").Value = 'AMT_CREDIT'
$ws.Range("B6").Value = 606
$ws.Range("A7").Value = 'YEARS_BIRTH'
$ws.Range("B7").Value = 558
$ws.Range("A8").Value = 'YEARS_EMPLOYED'
$ws.Range("B8").Value = 478
$ws.Range("A9").Value = 'OCCUPATION_TYPE'
$ws.Range("B9").Value = 457
$ws.Range("A10").Value = 'YEARS_REGISTRATION'
$ws.Range("B10").Value = 454
$ws.Range("A11").Value = 'AMT_ANNUITY'
$ws.Range("B11").Value = 427
$ws.Range("A12").Value = 'YEARS_ID_PUBLISH'
$ws.Range("B12").Value = 406
$ws.Range("A13").Value = 'REGION_POPULATION_RELATIVE'
$ws.Range("B13").Value = 345
$ws.Range("A14").Value = 'YEARS_LAST_PHONE_CHANGE'
$ws.Range("B14").Value = 296
$ws.Range("A15").Value = 'OWN_CAR_AGE'
$ws.Range("B15").Value = 269
$ws.Range("A16").Value = 'bur_cnt_active'
$ws.Range("B16").Value = 249
$ws.Range("A17").Value = 'HOUR_APPR_PROCESS_START'
$ws.Range("B17").Value = 215
$ws.Range("A18").Value = 'prev_cnt_consumer_approved'
$ws.Range("B18").Value = 213
$ws.Range("A19").Value = 'prev_cnt_cash_refused'
$ws.Range("B19").Value = 180
$ws.Range("A20").Value = 'CODE_GENDER'
$ws.Range("B20").Value = 113
$ws.Range("A21").Value = 'NAME_FAMILY_STATUS'
$ws.Range("B21").Value = 84
$ws.Range("A22").Value = 'prev_cnt_revolving_refused'
$ws.Range("B22").Value = 80
$ws.Range("A23").Value = 'REGION_RATING_CLIENT_W_CITY'
$ws.Range("B23").Value = 79
$ws.Range("A24").Value = 'NAME_EDUCATION_TYPE'
$ws.Range("B24").Value = 65
$ws.Range("A25").Value = 'FLAG_WORK_PHONE'
$ws.Range("B25").Value = 64
$ws.Range("A26").Value = 'NAME_CONTRACT_TYPE'
$ws.Range("B26").Value = 62
$ws.Range("A27").Value = 'NAME_HOUSING_TYPE'
$ws.Range("B27").Value = 61
$ws.Range("A28").Value = 'FLAG_DOCUMENT_3'
$ws.Range("B28").Value = 54
$ws.Range("A29").Value = 'REG_CITY_NOT_LIVE_CITY'
$ws.Range("B29").Value = 50
$ws.Range("A30").Value = 'bur_has_history'
$ws.Range("B30").Value = 43
$ws.Range("A31").Value = 'prev_cnt_revolving_canceled'
$ws.Range("B31").Value = 37
$ws.Range("A32").Value = 'prev_has_history'
$ws.Range("B32").Value = 36
$ws.Range("A33").Value = 'FLAG_PHONE'
$ws.Range("B33").Value = 34
$ws.Range("A34").Value = 'NAME_INCOME_TYPE'
$ws.Range("B34").Value = 28
$ws.Range("A35").Value = 'REG_CITY_NOT_WORK_CITY'
$ws.Range("B35").Value = 21
$ws.Range("A36").Value = 'FLAG_DOCUMENT_6'
$ws.Range("B36").Value = 8
$ws.Range("A37").Value = 'FLAG_OWN_CAR'
$ws.Range("B37").Value = 2
